# Add a new "Thank You" slide at the end of the deck, using the
# SECTION_HEADER custom layout (slideLayout2.xml), matching the layout
# already used for other section-style slides in this theme.

$p = $ppt.ActivePresentation

# SECTION_HEADER is the 2nd custom layout defined on the slide master.
$sectionHeaderLayout = $p.SlideMaster.CustomLayouts.Item(2)

# Append the new slide after the current last slide.
$newIndex = $p.Slides.Count + 1
$slide = $p.Slides.AddSlide($newIndex, $sectionHeaderLayout)

# Fill in the title placeholder text.
$title = $slide.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "THANKYOU"
$title.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# Materialize a notes page for the new slide (empty body, matching the
# other slides' generated notes pages).
$notes = $slide.NotesPage
$notesBody = $notes.Shapes.AddTextbox(1, 10, 10, 100, 100)
$notesBody.TextFrame.TextRange.Text = ""
